$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Öd1" and "Öd2" columns (B:C) entirely - remaining columns
# (Quiz, Vize, Fin, ORT) shift left into B:E.
$ws.Range("B:C").Delete()

# Delete the rows for students 220502025 and 220502015 (rows 3 and 4)
# - remaining rows shift up.
$ws.Range("3:4").Delete()
